# Chapter_4_Table_S4.7.xlsx - final supplementary update
# - title text: "cross validation" -> "cross-validation"
# - add a 5th table column "Hyperparamater 4" with per-model tree counts
# - update several hyperparameter values (LASSO/ENET lambda, NB adjust,
#   RF mtry/min.node.size, RF(ET) mtry)
# - widen column E and grow the data-row heights to fit the new column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Title text fix in A1: "cross validation" -> "cross-validation"
#    A1 holds rich text (bold "Supplementary Table S4.7:" + a regular
#    run for the rest). Re-apply both runs' formatting after the text
#    write so the bold/non-bold split is preserved.
# ---------------------------------------------------------------------
$titleRange = $ws.Range("A1")
$boldPart = "Supplementary Table S4.7:"
$restPart = " Selected hyperparamaters for models trained in the training set via five-fold cross-validation"
$newTitle = $boldPart + $restPart
$titleRange.Value = $newTitle

$boldChars = $titleRange.Characters(1, $boldPart.Length)
$boldChars.Font.Bold = $true

$restChars = $titleRange.Characters($boldPart.Length + 1, $restPart.Length)
$restChars.Font.Bold = $false

# ---------------------------------------------------------------------
# 2) Add the 5th table column ("Hyperparamater 4") to the ListObject.
#    Writing the header cell value keeps the table column name in sync.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null
$ws.Range("E2").Value = "Hyperparamater 4"

# ---------------------------------------------------------------------
# 3) Fill in column E (new "Hyperparamater 4") for every data row.
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("E8").Value = "num.trees = 751"
$ws.Range("E9").Value = "num.trees  = 3001"
$ws.Range("E10").Value = "NA"

# ---------------------------------------------------------------------
# 4) Update the changed hyperparameter values in columns B/C/D.
# ---------------------------------------------------------------------
$ws.Range("C5").Value = "lambda = 0.003810342"   # LASSO lambda
$ws.Range("C6").Value = "lambda = 0.05741517"    # ENET lambda
$ws.Range("D7").Value = "adjust = 0.3"           # NB adjust
$ws.Range("B8").Value = "mtry = 2"               # RF mtry
$ws.Range("D8").Value = "min.node.size = 3"      # RF min.node.size
$ws.Range("B9").Value = "mtry = 50"              # RF (ET) mtry

# ---------------------------------------------------------------------
# 5) Layout: widen the new column E and grow the data row heights so
#    the extra column's text fits.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 28

for ($r = 2; $r -le 10; $r++) {
    $ws.Rows.Item($r).RowHeight = 24.95
}
